$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Apoe"
$ws.Cells.Item(2,3).Value = "Sorl1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 100.8373843333333
$ws.Cells.Item(2,8).Value = 302.512153
$ws.Cells.Item(2,9).Value = 0.6551985585448407
$ws.Cells.Item(2,10).Value = 0.6551985585448408
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 86.09667233333334
$ws.Cells.Item(2,14).Value = 258.290017
$ws.Cells.Item(2,15).Value = 0.9236051571299395
$ws.Cells.Item(2,16).Value = 0.9236051571299394
$ws.Cells.Item(2,17).Value = 8681.763237897401
$ws.Cells.Item(2,18).Value = 78135.86914107662
$ws.Cells.Item(2,19).Value = 0.6051447676161175
$ws.Cells.Item(2,20).Value = 0.6051447676161175

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Apoe"
$ws.Cells.Item(3,3).Value = "Sorl1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 100.8373843333333
$ws.Cells.Item(3,8).Value = 302.512153
$ws.Cells.Item(3,9).Value = 0.6551985585448407
$ws.Cells.Item(3,10).Value = 0.6551985585448408
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 4.911922333333334
$ws.Cells.Item(3,14).Value = 14.735767
$ws.Cells.Item(3,15).Value = 0.05269282395635592
$ws.Cells.Item(3,16).Value = 0.05269282395635591
$ws.Cells.Item(3,17).Value = 495.3054001418168
$ws.Cells.Item(3,18).Value = 4457.748601276351
$ws.Cells.Item(3,19).Value = 0.03452426230186145
$ws.Cells.Item(3,20).Value = 0.03452426230186145

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Apoe"
$ws.Cells.Item(4,3).Value = "Sorl1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 100.8373843333333
$ws.Cells.Item(4,8).Value = 302.512153
$ws.Cells.Item(4,9).Value = 0.6551985585448407
$ws.Cells.Item(4,10).Value = 0.6551985585448408
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.209456
$ws.Cells.Item(4,14).Value = 6.628368
$ws.Cells.Item(4,15).Value = 0.02370201891370452
$ws.Cells.Item(4,16).Value = 0.02370201891370452
$ws.Cells.Item(4,17).Value = 222.7957638395893
$ws.Cells.Item(4,18).Value = 2005.161874556304
$ws.Cells.Item(4,19).Value = 0.01552952862686175
$ws.Cells.Item(4,20).Value = 0.01552952862686175

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Apoe"
$ws.Cells.Item(5,3).Value = "Sorl1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 23.90796933333333
$ws.Cells.Item(5,8).Value = 71.72390799999999
$ws.Cells.Item(5,9).Value = 0.1553438454249564
$ws.Cells.Item(5,10).Value = 0.1553438454249564
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 86.09667233333334
$ws.Cells.Item(5,14).Value = 258.290017
$ws.Cells.Item(5,15).Value = 0.9236051571299395
$ws.Cells.Item(5,16).Value = 0.9236051571299394
$ws.Cells.Item(5,17).Value = 2058.396601847382
$ws.Cells.Item(5,18).Value = 18525.56941662644
$ws.Cells.Item(5,19).Value = 0.1434763767628859
$ws.Cells.Item(5,20).Value = 0.1434763767628859

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Apoe"
$ws.Cells.Item(6,3).Value = "Sorl1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 23.90796933333333
$ws.Cells.Item(6,8).Value = 71.72390799999999
$ws.Cells.Item(6,9).Value = 0.1553438454249564
$ws.Cells.Item(6,10).Value = 0.1553438454249564
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 4.911922333333334
$ws.Cells.Item(6,14).Value = 14.735767
$ws.Cells.Item(6,15).Value = 0.05269282395635592
$ws.Cells.Item(6,16).Value = 0.05269282395635591
$ws.Cells.Item(6,17).Value = 117.4340885130485
$ws.Cells.Item(6,18).Value = 1056.906796617436
$ws.Cells.Item(6,19).Value = 0.008185505899680593
$ws.Cells.Item(6,20).Value = 0.008185505899680593

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Apoe"
$ws.Cells.Item(7,3).Value = "Sorl1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 23.90796933333333
$ws.Cells.Item(7,8).Value = 71.72390799999999
$ws.Cells.Item(7,9).Value = 0.1553438454249564
$ws.Cells.Item(7,10).Value = 0.1553438454249564
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.209456
$ws.Cells.Item(7,14).Value = 6.628368
$ws.Cells.Item(7,15).Value = 0.02370201891370452
$ws.Cells.Item(7,16).Value = 0.02370201891370452
$ws.Cells.Item(7,17).Value = 52.82360629134932
$ws.Cells.Item(7,18).Value = 475.412456622144
$ws.Cells.Item(7,19).Value = 0.003681962762389908
$ws.Cells.Item(7,20).Value = 0.003681962762389908

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Apoe"
$ws.Cells.Item(8,3).Value = "Sorl1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 29.15819666666667
$ws.Cells.Item(8,8).Value = 87.47459000000001
$ws.Cells.Item(8,9).Value = 0.1894575960302029
$ws.Cells.Item(8,10).Value = 0.1894575960302029
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 86.09667233333334
$ws.Cells.Item(8,14).Value = 258.290017
$ws.Cells.Item(8,15).Value = 0.9236051571299395
$ws.Cells.Item(8,16).Value = 0.9236051571299394
$ws.Cells.Item(8,17).Value = 2510.423704240893
$ws.Cells.Item(8,18).Value = 22593.81333816803
$ws.Cells.Item(8,19).Value = 0.1749840127509362
$ws.Cells.Item(8,20).Value = 0.1749840127509362

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Apoe"
$ws.Cells.Item(9,3).Value = "Sorl1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 29.15819666666667
$ws.Cells.Item(9,8).Value = 87.47459000000001
$ws.Cells.Item(9,9).Value = 0.1894575960302029
$ws.Cells.Item(9,10).Value = 0.1894575960302029
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 4.911922333333334
$ws.Cells.Item(9,14).Value = 14.735767
$ws.Cells.Item(9,15).Value = 0.05269282395635592
$ws.Cells.Item(9,16).Value = 0.05269282395635591
$ws.Cells.Item(9,17).Value = 143.2227974067256
$ws.Cells.Item(9,18).Value = 1289.00517666053
$ws.Cells.Item(9,19).Value = 0.009983055754813879
$ws.Cells.Item(9,20).Value = 0.009983055754813877

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Apoe"
$ws.Cells.Item(10,3).Value = "Sorl1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 29.15819666666667
$ws.Cells.Item(10,8).Value = 87.47459000000001
$ws.Cells.Item(10,9).Value = 0.1894575960302029
$ws.Cells.Item(10,10).Value = 0.1894575960302029
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.209456
$ws.Cells.Item(10,14).Value = 6.628368
$ws.Cells.Item(10,15).Value = 0.02370201891370452
$ws.Cells.Item(10,16).Value = 0.02370201891370452
$ws.Cells.Item(10,17).Value = 64.42375257434666
$ws.Cells.Item(10,18).Value = 579.8137731691201
$ws.Cells.Item(10,19).Value = 0.00449052752445286
$ws.Cells.Item(10,20).Value = 0.00449052752445286
